# Add "genre" / "set_in" columns to the "books" sheet and refresh the
# view-state (active sheet/cell, column width) to match the authored
# workbook state.

$wb = $excel.ActiveWorkbook

$books = $wb.Worksheets.Item("books")
$publisher = $wb.Worksheets.Item("publisher")

# --- New columns D (genre) and E (set_in) on the "books" sheet ----------

$books.Range("D1").Value = "genre"
$books.Range("E1").Value = "set_in"

# Albert Camus / The Plague
$books.Range("D2").Value = "Philosophical novel; absurdist novel"
$books.Range("E2").Value = "French Algeria"

# George Orwell / 1984
$books.Range("D3").Value = "dystopian novel;cautionary tale"
$books.Range("E3").Value = "United Kingdom"

# Match the font already used by the rest of the header/data cells so the
# new columns look consistent with the existing table.
$books.Range("D1:E3").Font.Name = "Arial"
$books.Range("D1:E3").Font.Size = 10

# --- View-state: widen the publisher name column -------------------------

$publisher.Columns.Item(2).ColumnWidth = 27.1666666667

# --- View-state: selections on each sheet --------------------------------

$publisher.Range("E29").Select()

# Switching to "books" makes it the active/tab-selected sheet, matching the
# authored workbook (publisher was active before the edit).
$books.Activate()
$books.Range("D31").Select()
